$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column C slightly to fit the new, longer test rows.
$ws.Columns.Item(3).ColumnWidth = 33.7

# --- New sprint-1 planning tests (breadth-first search / isFriend) ---
# Values are entered in the same left-to-right, row-by-row order the
# original author used, so shared strings line up with the source edit.

$ws.Range("B17").Value = "Knowledge basis"
$ws.Range("C17").Value = "breadth_search('JoseCid','Joao',L)"

$ws.Range("B18").Value = "Knowledge basis"
$ws.Range("C18").Value = "isFriend('Tiago','Tiago')."
$ws.Range("E18").Value = "False"
$ws.Range("D18").Value = "False ('Manuel' not found)"

$ws.Range("D17").Value = "L=['Joao','Diogo','Francisco','JoseCid']"
$ws.Range("E17").Value = "L=['Joao','Diogo','Francisco','JoseCid']"

$ws.Range("B19").Value = "Knowledge basis"
$ws.Range("C19").Value = "breadth_search('Artur',Stephanie',L)"

$ws.Range("D19").Value = "L=['Artur','Tiago','Stephanie']"
$ws.Range("E19").Value = "L=['Artur','Tiago','Stephanie']"

# New empty, underlined placeholder cell further down the sheet (mirrors C30).
$ws.Range("C32").Font.Underline = 2

# Move the active selection to the new placeholder cell, same as the author left it.
$ws.Range("C32").Select()
